$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows whose content was swapped/rotated in the source crawl ---

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value2 = "6346775"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value2 = "Tempo Premium 4-lagig 24 Rollen"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-premium-4-lagig-24-rollen/p/6346775"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "24Rol"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value2 = "Tempo"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value2 = "23.25"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value2 = "0.97/1Rol"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value2 = "Preis pro 1 Rolle"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value2 = "0.97"
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value2 = "1Rol"
$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N12").NumberFormat = "@"
$ws.Range("N12").Value2 = "Tempo Premium 4-lagig 24 Rollen 1 + 1 Aktion 23.25 Schweizer Franken"
$ws.Range("E12").Value2 = 7
$ws.Range("F12").Value2 = 5

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value2 = "6467117"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value2 = "Prix Garantie Toilettenpapier 3-lagig 10 Rollen"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/prix-garantie-toilettenpapier-3-lagig-10-rollen/p/6467117"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "10Rol"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value2 = "Coop"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value2 = "4.00"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value2 = "0.40/1Rol"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value2 = "Preis pro 1 Rolle"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value2 = "0.40"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value2 = "1Rol"
$ws.Range("M13").NumberFormat = "@"
$ws.Range("M13").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N13").NumberFormat = "@"
$ws.Range("N13").Value2 = "Prix Garantie Toilettenpapier 3-lagig 10 Rollen 4.00 Schweizer Franken"
$ws.Range("E13").Value2 = 18
$ws.Range("F13").Value2 = 4.5

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value2 = "4403154"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value2 = "Hakle Toilettenpapier Pflegende Sauberkeit 4-lagig 9 Rollen"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/hakle-toilettenpapier-pflegende-sauberkeit-4-lagig-9-rollen/p/4403154"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "9Rol"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value2 = "Hakle"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value2 = "8.90"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value2 = "0.99/1Rol"
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value2 = "Preis pro 1 Rolle"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value2 = "0.99"
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value2 = "1Rol"
$ws.Range("M22").NumberFormat = "@"
$ws.Range("M22").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N22").NumberFormat = "@"
$ws.Range("N22").Value2 = "Hakle Toilettenpapier Pflegende Sauberkeit 4-lagig 9 Rollen 8.90 Schweizer Franken"
$ws.Range("E22").Value2 = 52
$ws.Range("F22").Value2 = 4.5

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value2 = "5939634"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value2 = "Kleenex Box Collection"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value2 = "/de/kosmetik-gesundheit/gesichtspflege/kosmetiktuecher-watte/kosmetiktuecher/kleenex-box-collection/p/5939634"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "48BLT"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value2 = "Kleenex"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value2 = "2.20"
$ws.Range("I23").Value2 = ""
$ws.Range("J23").Value2 = ""
$ws.Range("K23").Value2 = ""
$ws.Range("L23").Value2 = ""
$ws.Range("M23").NumberFormat = "@"
$ws.Range("M23").Value2 = "['kosmetik-gesundheit', 'gesichtspflege', 'kosmetiktuecher-watte', 'kosmetiktuecher']"
$ws.Range("N23").NumberFormat = "@"
$ws.Range("N23").Value2 = "Kleenex Box Collection 2.20 Schweizer Franken"
$ws.Range("E23").Value2 = 11
$ws.Range("F23").Value2 = 4.5

# Row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value2 = "5985677"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value2 = "Tempo feucht Aloe sanft &amp; sensitiv 3x  42ST"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/tempo-feucht-aloe-sanft-sensitiv/p/5985677"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "3x 42ST"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value2 = "Tempo"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value2 = "6.90"
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value2 = "0.05/1ST"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value2 = "Preis pro 1 Stück"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value2 = "0.05"
$ws.Range("L24").NumberFormat = "@"
$ws.Range("L24").Value2 = "1ST"
$ws.Range("M24").NumberFormat = "@"
$ws.Range("M24").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N24").NumberFormat = "@"
$ws.Range("N24").Value2 = "Tempo feucht Aloe sanft &amp; sensitiv 3x  42ST 33% Aktion 6.90 Schweizer Franken statt 10.35 Schweizer Franken"
$ws.Range("E24").Value2 = 12
$ws.Range("F24").Value2 = 5

# Row 25
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value2 = "6833369"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value2 = "Tempo Taschentücher Box 100 Stück"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-taschentuecher-box-100-stueck/p/6833369"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "100ST"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value2 = "Tempo"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value2 = "3.95"
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value2 = "0.04/1ST"
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value2 = "Preis pro 1 Stück"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value2 = "0.04"
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value2 = "1ST"
$ws.Range("M25").NumberFormat = "@"
$ws.Range("M25").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N25").NumberFormat = "@"
$ws.Range("N25").Value2 = "Tempo Taschentücher Box 100 Stück 3.95 Schweizer Franken"
$ws.Range("E25").Value2 = 1
$ws.Range("F25").Value2 = 5

# Row 34
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value2 = "4489513"
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value2 = "Tempo Toilettenpapier white 3-lagig 32 Rollen"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-toilettenpapier-white-3-lagig-32-rollen/p/4489513"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "32Rol"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value2 = "Tempo"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value2 = "31.00"
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value2 = "0.97/1Rol"
$ws.Range("J34").NumberFormat = "@"
$ws.Range("J34").Value2 = "Preis pro 1 Rolle"
$ws.Range("K34").NumberFormat = "@"
$ws.Range("K34").Value2 = "0.97"
$ws.Range("L34").NumberFormat = "@"
$ws.Range("L34").Value2 = "1Rol"
$ws.Range("M34").NumberFormat = "@"
$ws.Range("M34").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N34").NumberFormat = "@"
$ws.Range("N34").Value2 = "Tempo Toilettenpapier white 3-lagig 32 Rollen 31.00 Schweizer Franken"
$ws.Range("E34").Value2 = 15
$ws.Range("F34").Value2 = 4

# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value2 = "6384992"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value2 = "Tempo Toilettenpapier Premium Kamille &amp; Aloe Vera 4-lagig 32 Rollen"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-toilettenpapier-premium-kamille-aloe-vera-4-lagig-32-rollen/p/6384992"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "32Rol"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value2 = "Tempo"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value2 = "31.00"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value2 = "0.97/1Rol"
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value2 = "Preis pro 1 Rolle"
$ws.Range("K35").NumberFormat = "@"
$ws.Range("K35").Value2 = "0.97"
$ws.Range("L35").NumberFormat = "@"
$ws.Range("L35").Value2 = "1Rol"
$ws.Range("M35").NumberFormat = "@"
$ws.Range("M35").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N35").NumberFormat = "@"
$ws.Range("N35").Value2 = "Tempo Toilettenpapier Premium Kamille &amp; Aloe Vera 4-lagig 32 Rollen 31.00 Schweizer Franken"
$ws.Range("E35").Value2 = 24
$ws.Range("F35").Value2 = 4

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value2 = "5985672"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value2 = "Tempo feucht sanft &amp; pur 3x  42ST"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/tempo-feucht-sanft-pur/p/5985672"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3x 42ST"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value2 = "Tempo"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value2 = "6.90"
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value2 = "0.05/1ST"
$ws.Range("J46").NumberFormat = "@"
$ws.Range("J46").Value2 = "Preis pro 1 Stück"
$ws.Range("K46").NumberFormat = "@"
$ws.Range("K46").Value2 = "0.05"
$ws.Range("L46").NumberFormat = "@"
$ws.Range("L46").Value2 = "1ST"
$ws.Range("M46").NumberFormat = "@"
$ws.Range("M46").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N46").NumberFormat = "@"
$ws.Range("N46").Value2 = "Tempo feucht sanft &amp; pur 3x  42ST 33% Aktion 6.90 Schweizer Franken statt 10.35 Schweizer Franken"
$ws.Range("E46").Value2 = 13
$ws.Range("F46").Value2 = 5

# Row 47
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value2 = "6286107"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value2 = "Oecoplan Haushaltpapier decor &amp; short, 2 Rollen 3-lagig"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/oecoplan-haushaltpapier-decor-short-2-rollen-3-lagig/p/6286107"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "200BLT"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value2 = "Coop"
$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value2 = "2.50"
$ws.Range("I47").Value2 = ""
$ws.Range("J47").Value2 = ""
$ws.Range("K47").Value2 = ""
$ws.Range("L47").Value2 = ""
$ws.Range("M47").NumberFormat = "@"
$ws.Range("M47").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N47").NumberFormat = "@"
$ws.Range("N47").Value2 = "Oecoplan Haushaltpapier decor &amp; short, 2 Rollen 3-lagig 2.50 Schweizer Franken"
$ws.Range("E47").Value2 = 14
$ws.Range("F47").Value2 = 3

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value2 = "6498160"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value2 = "subito Haushaltspapier weiss 8 Rollen"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/subito-haushaltspapier-weiss-8-rollen/p/6498160"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "400BLT"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value2 = "subito"
$ws.Range("H49").NumberFormat = "@"
$ws.Range("H49").Value2 = "10.95"
$ws.Range("I49").Value2 = ""
$ws.Range("J49").Value2 = ""
$ws.Range("K49").Value2 = ""
$ws.Range("L49").Value2 = ""
$ws.Range("M49").NumberFormat = "@"
$ws.Range("M49").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N49").NumberFormat = "@"
$ws.Range("N49").Value2 = "subito Haushaltspapier weiss 8 Rollen 10.95 Schweizer Franken"
$ws.Range("E49").Value2 = 3
$ws.Range("F49").Value2 = 4.5

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value2 = "3373453"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value2 = "Hakle feucht klassisch"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/hakle-feucht-klassisch/p/3373453"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "42ST"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value2 = "Hakle"
$ws.Range("H50").NumberFormat = "@"
$ws.Range("H50").Value2 = "3.20"
$ws.Range("I50").NumberFormat = "@"
$ws.Range("I50").Value2 = "0.08/1ST"
$ws.Range("J50").NumberFormat = "@"
$ws.Range("J50").Value2 = "Preis pro 1 Stück"
$ws.Range("K50").NumberFormat = "@"
$ws.Range("K50").Value2 = "0.08"
$ws.Range("L50").NumberFormat = "@"
$ws.Range("L50").Value2 = "1ST"
$ws.Range("M50").NumberFormat = "@"
$ws.Range("M50").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N50").NumberFormat = "@"
$ws.Range("N50").Value2 = "Hakle feucht klassisch 3.20 Schweizer Franken"
$ws.Range("E50").Value2 = 21
$ws.Range("F50").Value2 = 4.5

# Row 54
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value2 = "6995914"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value2 = "Tela Futura FSC 3-lagig 30 Rollen"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tela-futura-fsc-3-lagig-30-rollen/p/6995914"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value2 = "30Rol"
$ws.Range("G54").NumberFormat = "@"
$ws.Range("G54").Value2 = "Tela"
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value2 = "13.95"
$ws.Range("I54").NumberFormat = "@"
$ws.Range("I54").Value2 = "0.47/1Rol"
$ws.Range("J54").NumberFormat = "@"
$ws.Range("J54").Value2 = "Preis pro 1 Rolle"
$ws.Range("K54").NumberFormat = "@"
$ws.Range("K54").Value2 = "0.47"
$ws.Range("L54").NumberFormat = "@"
$ws.Range("L54").Value2 = "1Rol"
$ws.Range("M54").NumberFormat = "@"
$ws.Range("M54").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N54").NumberFormat = "@"
$ws.Range("N54").Value2 = "Tela Futura FSC 3-lagig 30 Rollen 50% Aktion 13.95 Schweizer Franken statt 28.00 Schweizer Franken"
$ws.Range("E54").Value2 = ""
$ws.Range("F54").Value2 = 0

# Row 55
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value2 = "6497243"
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value2 = "Super Soft WC-Papier Sensation 3-lagig 16 Rollen"
$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/super-soft-wc-papier-sensation-3-lagig-16-rollen/p/6497243"
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value2 = "16Rol"
$ws.Range("G55").NumberFormat = "@"
$ws.Range("G55").Value2 = "Super Soft"
$ws.Range("H55").NumberFormat = "@"
$ws.Range("H55").Value2 = "12.95"
$ws.Range("I55").NumberFormat = "@"
$ws.Range("I55").Value2 = "0.81/1Rol"
$ws.Range("J55").NumberFormat = "@"
$ws.Range("J55").Value2 = "Preis pro 1 Rolle"
$ws.Range("K55").NumberFormat = "@"
$ws.Range("K55").Value2 = "0.81"
$ws.Range("L55").NumberFormat = "@"
$ws.Range("L55").Value2 = "1Rol"
$ws.Range("M55").NumberFormat = "@"
$ws.Range("M55").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N55").NumberFormat = "@"
$ws.Range("N55").Value2 = "Super Soft WC-Papier Sensation 3-lagig 16 Rollen 12.95 Schweizer Franken"
$ws.Range("E55").Value2 = 4
$ws.Range("F55").Value2 = 3.5

# Row 56
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value2 = "6283679"
$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value2 = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen"
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-toilettenpapier-camomille-weiss-4-lagig-6-rollen/p/6283679"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value2 = "6Rol"
$ws.Range("G56").NumberFormat = "@"
$ws.Range("G56").Value2 = "Coop"
$ws.Range("H56").NumberFormat = "@"
$ws.Range("H56").Value2 = "4.20"
$ws.Range("I56").NumberFormat = "@"
$ws.Range("I56").Value2 = "0.70/1Rol"
$ws.Range("J56").NumberFormat = "@"
$ws.Range("J56").Value2 = "Preis pro 1 Rolle"
$ws.Range("K56").NumberFormat = "@"
$ws.Range("K56").Value2 = "0.70"
$ws.Range("L56").NumberFormat = "@"
$ws.Range("L56").Value2 = "1Rol"
$ws.Range("M56").NumberFormat = "@"
$ws.Range("M56").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N56").NumberFormat = "@"
$ws.Range("N56").Value2 = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen 4.20 Schweizer Franken"
$ws.Range("E56").Value2 = 13
$ws.Range("F56").Value2 = 4

# Row 59
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value2 = "4947421"
$ws.Range("B59").NumberFormat = "@"
$ws.Range("B59").Value2 = "Oecoplan Taschentuch Calendula Box"
$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value2 = "/de/inspiration-geschenke/saisonale-promotionen/gesundheit/oecoplan-taschentuch-calendula-box/p/4947421"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value2 = "80ST"
$ws.Range("G59").NumberFormat = "@"
$ws.Range("G59").Value2 = "Coop"
$ws.Range("H59").NumberFormat = "@"
$ws.Range("H59").Value2 = "2.30"
$ws.Range("I59").NumberFormat = "@"
$ws.Range("I59").Value2 = "0.03/1ST"
$ws.Range("J59").NumberFormat = "@"
$ws.Range("J59").Value2 = "Preis pro 1 Stück"
$ws.Range("K59").NumberFormat = "@"
$ws.Range("K59").Value2 = "0.03"
$ws.Range("L59").NumberFormat = "@"
$ws.Range("L59").Value2 = "1ST"
$ws.Range("M59").NumberFormat = "@"
$ws.Range("M59").Value2 = "['inspiration-geschenke', 'saisonale-promotionen', 'gesundheit']"
$ws.Range("N59").NumberFormat = "@"
$ws.Range("N59").Value2 = "Oecoplan Taschentuch Calendula Box 2.30 Schweizer Franken"
$ws.Range("E59").Value2 = 17
$ws.Range("F59").Value2 = 4

# Row 60
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value2 = "6497242"
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value2 = "Super Soft WC-Papier Sensation, 4 Rollen 3-lagig"
$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/super-soft-wc-papier-sensation-4-rollen-3-lagig/p/6497242"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value2 = "4Rol"
$ws.Range("G60").NumberFormat = "@"
$ws.Range("G60").Value2 = "Super Soft"
$ws.Range("H60").NumberFormat = "@"
$ws.Range("H60").Value2 = "3.25"
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value2 = "0.81/1Rol"
$ws.Range("J60").NumberFormat = "@"
$ws.Range("J60").Value2 = "Preis pro 1 Rolle"
$ws.Range("K60").NumberFormat = "@"
$ws.Range("K60").Value2 = "0.81"
$ws.Range("L60").NumberFormat = "@"
$ws.Range("L60").Value2 = "1Rol"
$ws.Range("M60").NumberFormat = "@"
$ws.Range("M60").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N60").NumberFormat = "@"
$ws.Range("N60").Value2 = "Super Soft WC-Papier Sensation, 4 Rollen 3-lagig 3.25 Schweizer Franken"
$ws.Range("E60").Value2 = 6
$ws.Range("F60").Value2 = 4.5

# Row 68
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value2 = "4489512"
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value2 = "Tempo Toilettenpapier weiss 3-lagig 24 Rollen"
$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tempo-toilettenpapier-weiss-3-lagig-24-rollen/p/4489512"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value2 = "24Rol"
$ws.Range("G68").NumberFormat = "@"
$ws.Range("G68").Value2 = "Tempo"
$ws.Range("H68").NumberFormat = "@"
$ws.Range("H68").Value2 = "23.25"
$ws.Range("I68").NumberFormat = "@"
$ws.Range("I68").Value2 = "0.97/1Rol"
$ws.Range("J68").NumberFormat = "@"
$ws.Range("J68").Value2 = "Preis pro 1 Rolle"
$ws.Range("K68").NumberFormat = "@"
$ws.Range("K68").Value2 = "0.97"
$ws.Range("L68").NumberFormat = "@"
$ws.Range("L68").Value2 = "1Rol"
$ws.Range("M68").NumberFormat = "@"
$ws.Range("M68").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N68").NumberFormat = "@"
$ws.Range("N68").Value2 = "Tempo Toilettenpapier weiss 3-lagig 24 Rollen 1+1 Aktion 23.25 Schweizer Franken"
$ws.Range("E68").Value2 = 15
$ws.Range("F68").Value2 = 3.5

# Row 69
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value2 = "6711017"
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value2 = "Tempo feuchte Limited Edition 3x  42ST"
$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/tempo-feuchte-limited-edition/p/6711017"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value2 = "3x 42ST"
$ws.Range("G69").NumberFormat = "@"
$ws.Range("G69").Value2 = "Tempo"
$ws.Range("H69").NumberFormat = "@"
$ws.Range("H69").Value2 = "6.90"
$ws.Range("I69").NumberFormat = "@"
$ws.Range("I69").Value2 = "0.05/1ST"
$ws.Range("J69").NumberFormat = "@"
$ws.Range("J69").Value2 = "Preis pro 1 Stück"
$ws.Range("K69").NumberFormat = "@"
$ws.Range("K69").Value2 = "0.05"
$ws.Range("L69").NumberFormat = "@"
$ws.Range("L69").Value2 = "1ST"
$ws.Range("M69").NumberFormat = "@"
$ws.Range("M69").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N69").NumberFormat = "@"
$ws.Range("N69").Value2 = "Tempo feuchte Limited Edition 3x  42ST 33% Aktion 6.90 Schweizer Franken statt 10.35 Schweizer Franken"
$ws.Range("E69").Value2 = 2
$ws.Range("F69").Value2 = 2.5

# Row 75
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value2 = "6498157"
$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value2 = "subito Haushaltspapier weiss 2 Rollen"
$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value2 = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/subito-haushaltspapier-weiss-2-rollen/p/6498157"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value2 = "100BLT"
$ws.Range("G75").NumberFormat = "@"
$ws.Range("G75").Value2 = "subito"
$ws.Range("H75").NumberFormat = "@"
$ws.Range("H75").Value2 = "2.75"
$ws.Range("I75").Value2 = ""
$ws.Range("J75").Value2 = ""
$ws.Range("K75").Value2 = ""
$ws.Range("L75").Value2 = ""
$ws.Range("M75").NumberFormat = "@"
$ws.Range("M75").Value2 = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N75").NumberFormat = "@"
$ws.Range("N75").Value2 = "subito Haushaltspapier weiss 2 Rollen 2.75 Schweizer Franken"
$ws.Range("E75").Value2 = ""
$ws.Range("F75").Value2 = 0

# Row 76
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value2 = "3640534"
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value2 = "Kandoo Feuchttücher Sensitive 50 Stück"
$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value2 = "/de/kosmetik-gesundheit/baby-kind/pflege-accessoires/feuchttuecher/kandoo-feuchttuecher-sensitive-50-stueck/p/3640534"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value2 = "50ST"
$ws.Range("G76").NumberFormat = "@"
$ws.Range("G76").Value2 = "Kandoo"
$ws.Range("H76").NumberFormat = "@"
$ws.Range("H76").Value2 = "3.75"
$ws.Range("I76").NumberFormat = "@"
$ws.Range("I76").Value2 = "0.08/1ST"
$ws.Range("J76").NumberFormat = "@"
$ws.Range("J76").Value2 = "Preis pro 1 Stück"
$ws.Range("K76").NumberFormat = "@"
$ws.Range("K76").Value2 = "0.08"
$ws.Range("L76").NumberFormat = "@"
$ws.Range("L76").Value2 = "1ST"
$ws.Range("M76").NumberFormat = "@"
$ws.Range("M76").Value2 = "['kosmetik-gesundheit', 'baby-kind', 'pflege-accessoires', 'feuchttuecher']"
$ws.Range("N76").NumberFormat = "@"
$ws.Range("N76").Value2 = "Kandoo Feuchttücher Sensitive 50 Stück 3.75 Schweizer Franken"
$ws.Range("E76").Value2 = 10
$ws.Range("F76").Value2 = 5

# --- Update crawl timestamp for every data row ---
$ws.Range("O2:O91").NumberFormat = "@"
$ws.Range("O2:O91").Value2 = "2022-07-16 20:56:15"
